$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "38.035.00"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "2.039.85"
$ws.Range("E3").Value = "  -0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'228.81"
$ws.Range("E5").Value = "  -0.18%  "

# Row 6
$ws.Range("E6").Value = "  -0.86%  "

# Row 7
$ws.Range("D7").Value = "'60.69"
$ws.Range("E7").Value = "  +3.41%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -0.75%  "

# Row 10
$ws.Range("D10").Value = "'0.0819"
$ws.Range("E10").Value = "  +1.23%  "

# Row 11
$ws.Range("E11").Value = "  +0.37%  "

# Row 12
$ws.Range("E12").Value = "  +0.01%  "

# Row 13
$ws.Range("D13").Value = "2.342.18"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("E14").Value = "  +2.83%  "

# Row 15
$ws.Range("E15").Value = "  +1.89%  "

# Row 16
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  -1.68%  "

# Row 17
$ws.Range("D17").Value = "2.025.56"
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").Value = "37.936.56"
$ws.Range("E18").Value = "  +0.10%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'5.99"
$ws.Range("E19").Value = "  -5.01%  "

# Row 20
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D20").Value = "'70.01"
$ws.Range("E20").Value = "  +0.55%  "

# Row 21
$ws.Range("E21").Value = "  -1.03%  "

# Row 22
$ws.Range("D22").Value = "'224.92"
$ws.Range("E22").Value = "  +0.15%  "

# Row 23
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.10%  "

# Row 24
$ws.Range("E24").Value = "  -0.19%  "

# Row 25
$ws.Range("D25").Value = "'2.25"
$ws.Range("E25").Value = "  +0.86%  "

# Row 26
$ws.Range("D26").Value = "'167.40"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("D28").Value = "'0.130"
$ws.Range("E28").Value = "  -2.99%  "

# Row 29
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
$ws.Range("D30").Value = "'1.29"
$ws.Range("E30").Value = "  -2.92%  "

# Row 31
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "  +0.93%  "

# Row 32
$ws.Range("D32").Value = "'2.17"
$ws.Range("E32").Value = "  +5.42%  "

# Row 33
$ws.Range("E33").Value = "  -2.04%  "

# Row 34
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.54"
$ws.Range("E34").Value = "  -0.87%  "

# Row 35
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0609"
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("D36").Value = "'6.47"
$ws.Range("E36").Value = "  +6.26%  "

# Row 37
$ws.Range("E37").Value = "  -1.76%  "

# Row 38
$ws.Range("E38").Value = "  +0.23%  "

# Row 39
$ws.Range("E39").Value = "  +0.06%  "

# Row 40
$ws.Range("D40").Value = "1.527.24"
$ws.Range("E40").Value = "  +2.26%  "

# Row 41
$ws.Range("D41").Value = "'17.45"

# Row 42
$ws.Range("D42").Value = "'0.0218"
$ws.Range("E42").Value = "  +0.46%  "

# Row 43
$ws.Range("D43").Value = "'96.57"
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("E44").Value = "  -1.73%  "

# Row 45
$ws.Range("D45").Value = "'0.0918"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46
$ws.Range("E46").Value = "  -2.82%  "

# Row 47
$ws.Range("D47").Value = "'3.98"
$ws.Range("E47").Value = "  -3.95%  "

# Row 48
$ws.Range("E48").Value = "  -0.47%  "

# Row 49
$ws.Range("E49").Value = "  +0.09%  "

# Row 50
$ws.Range("E50").Value = "  +0.47%  "

# Row 51
$ws.Range("D51").Value = "2.230.43"
$ws.Range("E51").Value = "  -0.53%  "
